$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.776.59'
$ws.Range('E2').Value = '  +0.89%  '

$ws.Range('D3').Value = '3.143.69'
$ws.Range('E3').Value = '  +1.08%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('E5').Value = '  +0.44%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.76'
$ws.Range('E6').Value = '  +0.48%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = '3.142.05'
$ws.Range('E8').Value = '  +1.25%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -0.30%  '

$ws.Range('E10').Value = '  +5.93%  '

$ws.Range('E11').Value = '  -0.36%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -2.01%  '

$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.16'
$ws.Range('E14').Value = '  +4.80%  '

$ws.Range('D15').Value = '3.664.13'
$ws.Range('E15').Value = '  +1.10%  '

$ws.Range('E16').Value = '  -1.33%  '

$ws.Range('D17').Value = '3.144.68'
$ws.Range('E17').Value = '  +1.21%  '

$ws.Range('D18').Value = '63.599.09'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.09'
$ws.Range('E19').Value = '  -0.87%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '463.46'
$ws.Range('E20').Value = '  -0.66%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.27'
$ws.Range('E21').Value = '  +0.78%  '

$ws.Range('E22').Value = '  +0.67%  '

$ws.Range('E23').Value = '  -1.05%  '

$ws.Range('E24').Value = '  -2.25%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.31'
$ws.Range('E25').Value = '  -0.82%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('E26').Value = '  +1.54%  '

$ws.Range('E27').Value = '  -0.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.18'
$ws.Range('E28').Value = '  +6.98%  '

$ws.Range('E29').Value = '  +0.36%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  -0.23%  '

$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.96'
$ws.Range('E32').Value = '  +1.43%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.01'
$ws.Range('E33').Value = '  +0.15%  '

$ws.Range('E34').Value = '  +0.48%  '

$ws.Range('D35').Value = '0.0₃0852'
$ws.Range('E35').Value = '  -1.96%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  -0.34%  '

$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.30'
$ws.Range('E37').Value = '  -4.60%  '

$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.35'
$ws.Range('E38').Value = '  +0.07%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.02'
$ws.Range('E39').Value = '  -0.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.95'
$ws.Range('E40').Value = '  +0.86%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '439.57'
$ws.Range('E41').Value = '  +0.95%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.77'
$ws.Range('E42').Value = '  +0.78%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0372'
$ws.Range('E43').Value = '  +0.89%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.914.83'
$ws.Range('E44').Value = '  +0.08%  '

$ws.Range('E45').Value = '  -0.82%  '

$ws.Range('E46').Value = '  -1.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.87'
$ws.Range('E47').Value = '  +5.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.11'
$ws.Range('E48').Value = '  +2.53%  '

$ws.Range('E49').Value = '  +0.00%  '

$ws.Range('E50').Value = '  -0.63%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.39'
$ws.Range('E51').Value = '  -0.86%  '
